# Update MarketoData sheet: refresh metric values and append new Marketo
# asset-type rows (Snippets .. Target Account Management).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B updates for existing rows (2-12) -----------------------------
$ws.Range("B2").Value  = 14
$ws.Range("B3").Value  = 10
$ws.Range("B4").Value  = 11
$ws.Range("B9").Value  = 6
$ws.Range("B10").Value = 6
$ws.Range("B11").Value = 0

# --- New rows (13-27) -------------------------------------------------------
$ws.Range("A13").Value = "Snippets"
$ws.Range("B13").Value = 8

$ws.Range("A14").Value = "Account Name"
$ws.Range("B14").Value = "send_box"

$ws.Range("A15").Value = "Models"
$ws.Range("B15").Value = 5

$ws.Range("A16").Value = "Leads"
$ws.Range("B16").Value = 7

$ws.Range("A17").Value = "Tokens"
$ws.Range("B17").Value = 5

$ws.Range("A18").Value = "Library"
$ws.Range("B18").Value = 4

$ws.Range("A19").Value = "Event Programs"
$ws.Range("B19").Value = 1

$ws.Range("A20").Value = "Nurture campaigns"
$ws.Range("B20").Value = 0

$ws.Range("A21").Value = "Segment Data"
$ws.Range("B21").Value = 10

$ws.Range("A22").Value = "Integration Data"
$ws.Range("B22").Value = 6

$ws.Range("A23").Value = "Interesting Moment_subscription"
$ws.Range("B23").Value = $true

$ws.Range("A24").Value = "Web Personalize"
$ws.Range("B24").Value = $true

$ws.Range("A25").Value = "All Batch Campaigns"
$ws.Range("B25").Value = 6

$ws.Range("A26").Value = "Images and Files"
$ws.Range("B26").Value = 4

$ws.Range("A27").Value = "Target Account Management"
$ws.Range("B27").Value = $true

# --- Formatting: row 20 gets a slightly larger, dark-grey font + taller row -
$f = $ws.Range("A20").Font
$f.Size = 12
$f.Color = 1250067
$ws.Rows.Item(20).RowHeight = 15.5

# --- View state: scroll down and select the last-touched cell --------------
$ws.Range("B26").Select()
